$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND: $find"
    }
    return $ok
}

# 1. Update the document date (appears near the top, "Date" style paragraph)
Replace-Text "2025-05-15" "2025-06-05"

# 2. Oxford comma: "empirical, experimental and simulation results"
Replace-Text "empirical, experimental and simulation results" "empirical, experimental, and simulation results"

# 3. "e.g. its DOI" -> "e.g., its DOI"
Replace-Text "permanent identifier (e.g. its DOI)" "permanent identifier (e.g., its DOI)"

# 4. "e.g. setting a seed" -> "e.g., setting a seed"
Replace-Text "appropriate means (e.g. setting a seed for the random number generator)" "appropriate means (e.g., setting a seed for the random number generator)"

# 5. "(e.g. estimation, ...)" and "(i.e. tables, ...)" in the same paragraph
Replace-Text "computational results (e.g. estimation, simulation, model solution, visualization, etc.)" "computational results (e.g., estimation, simulation, model solution, visualization, etc.)"
Replace-Text "approved online appendices (i.e. tables, figures, in-text numbers)" "approved online appendices (i.e., tables, figures, in-text numbers)"

# 6. "reasonable timeframe" -> "reasonable time frame"
Replace-Text "reasonable timeframe" "reasonable time frame"

# 7. "strongly encouraged, and must be provided ... (e.g. to limit ...)" and "(e.g. tables and figures)"
Replace-Text "strongly encouraged, and must be provided at the Data Editor" "strongly encouraged and must be provided at the Data Editor"
Replace-Text "specific request (e.g. to limit the number of human intervention steps)" "specific request (e.g., to limit the number of human intervention steps)"
Replace-Text "save all exhibits (e.g. tables and figures)" "save all exhibits (e.g., tables and figures)"

# 8. "(such as, for example, precise specification...)" and "(e.g., ArcGIS)" -> "(e.g., ArcGIS or MS Excel)"
Replace-Text "should be taken (such as, for example, precise specification" "should be taken (such as, e.g., precise specification"
Replace-Text "using scripts (e.g., ArcGIS) is discouraged" "using scripts (e.g., ArcGIS or MS Excel) is discouraged"

# 9. "Authors who collect primary data (e.g. via experiment or survey)"
Replace-Text "Authors who collect primary data (e.g. via experiment or survey)" "Authors who collect primary data (e.g., via experiment or survey)"

# 10. Delete the entire paragraph beginning "Other repositories and archives may be acceptable..."
$rng = $d.Content
$found = $rng.Find.Execute("Other repositories and archives may be acceptable", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(1)
    [void]$rng.Expand(4)
    $rng.Delete()
} else {
    Write-Host "NOT FOUND: paragraph to delete"
}

# 11. "In cases where data cannot be published in an openly accessible trusted data repository," -> add "like the JPE dataverse"
Replace-Text "openly accessible trusted data repository, authors who have requested an exemption" "openly accessible trusted data repository like the JPE dataverse, authors who have requested an exemption"
